# Weekly driver report update for 2025-04-20
# Updates the "Bad Drivers" table (rows 3-4 + totals row 5) and the
# "Good Drivers" table (rows 13-21) on the active "Driver Summary" sheet
# to reflect the refreshed roaming-impact source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Bad Drivers table
# ---------------------------------------------------------------------

# Row 3
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.10.1"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 1805
$ws.Range("D3").Value = 93.09999999999999

# Row 4
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.160.0.3"
$ws.Range("B4").Value = 46
$ws.Range("C4").Value = 7740
$ws.Range("D4").Value = 97.2

# Row 5 - Totals (client count total is unchanged, only critical minutes grows)
$ws.Range("C5").Value = 9545

# ---------------------------------------------------------------------
# Good Drivers table
#
# The E column ("Driver Vintage") holds its dates as literal text, not
# real Excel date serials, in this workbook. Pre-formatting the cells as
# Text before assigning the date-looking strings keeps Excel from
# auto-converting them into date serial numbers.
# ---------------------------------------------------------------------
$ws.Range("E13:E21").NumberFormat = "@"

# Row 13
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B13").Value = 445055
$ws.Range("D13").Value = 99.90000000000001
$ws.Range("E13").Value = "2024-11-10"

# Row 14
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B14").Value = 10661
$ws.Range("E14").Value = "2022-08-29"

# Row 15
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B15").Value = 14239
$ws.Range("E15").Value = "2022-05-23"

# Row 16
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B16").Value = 265400
$ws.Range("E16").Value = "2022-05-01"

# Row 17
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B17").Value = 77849
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").Value = "2021-08-18"

# Row 18
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B18").Value = 34244
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "2021-04-27"

# Row 19
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B19").Value = 59673
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "2020-08-05"

# Row 20
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B20").Value = 113652
$ws.Range("E20").Value = "2020-01-06"

# Row 21
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B21").Value = 56018
$ws.Range("E21").Value = "2019-12-14"
